# Apply the Feb 2024 绩效表 (performance sheet) updates:
#  - Fill in K7/K8/K9 绩效系数 (performance coefficients)
#  - Shift the C/E/G task rows 10/11/12 down one slot and add a new
#    "雅威1月业务审核" entry, filling their K coefficients
#  - Correct D13's business type and fill its K coefficient
#  - Update the 组员1/其他业务 summary row formulas & totals (row 29/30)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("绩效表")

# --- K column 绩效系数 for rows 7-13 -------------------------------------
$ws.Range("K7").Value = 0.15
$ws.Range("K8").Value = 0.083
$ws.Range("K9").Value = 0.25

# --- Row 10: order number / note / title updated --------------------------
$ws.Range("C10").Value = "A2023112804"
$ws.Range("E10").Value = "sci3-5"
$ws.Range("G10").Value = "质谱+网络药理学分析"
$ws.Range("K10").Value = 0.15

# --- Row 11: order number / note / title updated --------------------------
$ws.Range("C11").Value = "周宇查询学者发文和 H 指数"
$ws.Range("E11").Value = ""
$ws.Range("G11").Value = "查询学者发文和 H 指数"
$ws.Range("K11").Value = 0.15

# --- Row 12: new task replacing the old one --------------------------------
$ws.Range("C12").Value = "雅威1月业务审核"
$ws.Range("G12").Value = "雅威1月业务审核"
$ws.Range("K12").Value = 0.008

# --- Row 13: business type correction + coefficient ------------------------
$ws.Range("D13").Value = "固定业务"
$ws.Range("K13").Value = 0.25

# --- Summary rows 29 (基本业务) / 30 (其他业务) -----------------------------
$ws.Range("F29").Value = 1
$ws.Range("H29").Value = 0.25
$ws.Range("I29").Value = "0.25=0.25"
$ws.Range("J29").Value = 1.041

$ws.Range("F30").Value = 6
$ws.Range("G30").Value = ""
$ws.Range("H30").Value = 0.791
$ws.Range("I30").Value = "0.15+0.083+0.25+0.15+0.15+0.008=0.791"
